# Update "想去人数" (column F) counts across all four sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 58
$ws1.Range("F3").Value = 58
$ws1.Range("F4").Value = 953
$ws1.Range("F5").Value = 1248
$ws1.Range("F6").Value = 1717
$ws1.Range("F9").Value = 2511
$ws1.Range("F10").Value = 727
$ws1.Range("F13").Value = 23
$ws1.Range("F16").Value = 227
$ws1.Range("F19").Value = 1229
$ws1.Range("F20").Value = 704
$ws1.Range("F21").Value = 3
$ws1.Range("F27").Value = 309
$ws1.Range("F28").Value = 309
$ws1.Range("F30").Value = 1771
$ws1.Range("F32").Value = 529
$ws1.Range("F33").Value = 520
$ws1.Range("F36").Value = 4564
$ws1.Range("F37").Value = 134

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 407
$ws2.Range("F8").Value = 63
$ws2.Range("F26").Value = 250

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1419
$ws3.Range("F6").Value = 518
$ws3.Range("F7").Value = 74
$ws3.Range("F8").Value = 189

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1419
$ws4.Range("F6").Value = 58
$ws4.Range("F8").Value = 58
$ws4.Range("F9").Value = 1248
$ws4.Range("F10").Value = 1717
$ws4.Range("F12").Value = 63
$ws4.Range("F16").Value = 2511
$ws4.Range("F17").Value = 727
$ws4.Range("F20").Value = 23
$ws4.Range("F24").Value = 227
$ws4.Range("F32").Value = 3
$ws4.Range("F36").Value = 189
$ws4.Range("F40").Value = 309
$ws4.Range("F41").Value = 1771
$ws4.Range("F42").Value = 250
$ws4.Range("F43").Value = 520
$ws4.Range("F46").Value = 4564
$ws4.Range("F47").Value = 134
